# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.847.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.421.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.799.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.413.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.828'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.691.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0917'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.39%  '
$ws.Range("E24").Value = '  +3.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("E32").Value = '  +15.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.78%  '
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0762'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.39%  '
$ws.Range("E37").Value = '  +2.85%  '
$ws.Range("E38").Value = '  +3.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '126.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("E42").Value = '  -3.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("E44").Value = '  +3.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.938.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("E49").Value = '  +16.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.13%  '
